$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original summary table lives in A1:D8. The sheet is extended with two
# more copies of that same table, pasted below (rows 11-18 and 21-28), each
# with the data cells (B:C for the data rows) cleared out, keeping only the
# labels/headers and formatting intact.

function Copy-Table([string]$destTopLeft) {
    $destCol = $destTopLeft.Substring(0,1)
    $destRow = [int]($destTopLeft.Substring(1))

    # values + formulas
    $ws.Range("A1:D8").Copy()
    $ws.Range($destTopLeft).PasteSpecial(-4104)

    # formats for the A:C block
    $ws.Range("A1:C8").Copy()
    $ws.Range($destTopLeft).PasteSpecial(-4122)

    # format for D1 (top merged header row)
    $ws.Range("D1").Copy()
    $ws.Range("D" + $destRow).PasteSpecial(-4122)

    # format for D6:D7 (the two "customFormat" data rows)
    $ws.Range("D6:D7").Copy()
    $ws.Range("D" + ($destRow + 5) + ":D" + ($destRow + 6)).PasteSpecial(-4122)

    # clear the data values (columns B:C, the 5 data rows) but keep formatting
    $ws.Range("B" + ($destRow + 3) + ":C" + ($destRow + 7)).ClearContents()
}

Copy-Table("A11")
Copy-Table("A21")

$excel.CutCopyMode = $false

$ws.Range("D26").Select()
